$wb = $excel.ActiveWorkbook

# "Yahoo" worksheet: remove the "4444" segment row (old row 3), which
# shifts the 1111/2222/3333 rows up by one, and append a new "5555"
# segment row at the bottom.
$ws = $wb.Worksheets.Item("Yahoo")
$ws.Rows.Item(3).Delete()

$ws.Range("A6").Value = 5555
$ws.Range("B6").Value = "This is simply a test"
$ws.Range("C6").Value = "Test 5555"

# Make "Yahoo" the active sheet/tab and move the selection to A7 (just
# below the newly added row), matching the saved view state.
$ws.Activate()
$ws.Range("A7").Select()
